# WTREGEN.xlsx update
# - "Data" sheet: append the newest weekly observation (row 96).
# - "SeriesInfo" sheet: refresh the FRED metadata (realtime dates,
#   observation_end, last_updated timestamp and popularity score).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Data" - add the new observation row
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Copy row 95's formatting down to the new row 96 first (keeps the same
# date number-format/border/font the rest of column A uses), then set
# the actual values for the new observation.
$wsData.Range("A95").Copy($wsData.Range("A96"))
$wsData.Range("A96").Value = 45133
$wsData.Range("B96").Value = 543.597

# ---------------------------------------------------------------------
# Sheet 2: "SeriesInfo" - refresh metadata fields
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# realtime_start / realtime_end / observation_end hold plain ISO date
# strings (not real dates) in the original file. Assigning a
# date-looking literal directly to a cell lets the engine auto-coerce
# it into a real date serial number, which we don't want. Instead,
# build the text in a scratch cell via a formula (a formula result is
# always kept as text/whatever type it evaluates to, never re-parsed
# as a date) and copy that value+type over to the destination - this
# avoids touching the destination cell's style too.
$helper = $wsInfo.Range("D1")

$helper.Formula = "=""2023-08-03"""
$helper.Copy($wsInfo.Range("B3"))

$helper.Formula = "=""2023-08-03"""
$helper.Copy($wsInfo.Range("B4"))

$helper.Formula = "=""2023-07-26"""
$helper.Copy($wsInfo.Range("B7"))

# Remove the scratch cell again so it doesn't linger in the sheet.
$helper.Value = ""

# last_updated already contains a UTC-offset suffix, so Excel leaves it
# as text without any extra coercion guard needed.
$wsInfo.Range("B14").Value = "2023-07-27 15:34:02-05"

# popularity is numeric.
$wsInfo.Range("B15").Value = 82
